$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Brand" header in column E, matching the style of the other headers
$ws.Range("E1").Value = "Brand"
$ws.Range("E1").Style = $ws.Range("A1").Style

# Brand values for existing rows 2-8 (alternating Bonapapa / Candyland)
$ws.Range("E2").Value = "Bonapapa"
$ws.Range("E3").Value = "Candyland"
$ws.Range("E4").Value = "Bonapapa"
$ws.Range("E5").Value = "Candyland"
$ws.Range("E6").Value = "Bonapapa"
$ws.Range("E7").Value = "Candyland"
$ws.Range("E8").Value = "Bonapapa"

# New row 9: a new shopkeeper record
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Ahmed"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "53454634634"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = "Candyland"
